$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "pdf_online" row's name/description from "pdf online" to "pdf online test"
$ws.Range("C4").Value = "pdf online test"
$ws.Range("D4").Value = "pdf online test"

# Update the saved selection/active cell shown when the sheet is reopened
$ws.Range("C10").Select()
